# Add two new columns (I: "I0", J: "IF") to the sheet, matching the
# header style already used by the other header cells (B1:H1), and
# fill in the data rows 2 and 3 for these new columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting of the existing header cell H1 (bold, centered,
# bordered style) onto the two new header cells so they match the rest
# of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data values for rows 2 and 3
$ws.Range("I2").Value = 3
$ws.Range("J2").Value = 4
$ws.Range("I3").Value = 7
$ws.Range("J3").Value = 7
